# Update handback/handoff timestamps to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for edb9c191-... row
$wsOverview.Range("G2").Value = "2016-09-07 11:32:39"

# zh-cn: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for edb9c191-... row
$wsZhCn.Range("H2").Value = "2016-09-07 11:32:33"
$wsZhCn.Range("K2").Value = "2016-09-07 11:33:01"

# de-de: "Correspond Handoff Datetime" for edb9c191-... row
$wsDeDe.Range("H2").Value = "2016-09-07 11:32:39"
# de-de: "Correspond Handback DateTime" for edb9c191-... row
$wsDeDe.Range("K2").Value = "2016-09-07 11:33:24"
